$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stores")

# New section / store rows added to support the sections combobox
$data = @(
    @("2", "1ST FLOOR", "q",  "q"),
    @("3", "1ST FLOOR", "a",  "a"),
    @("4", "FOOD",      "qw", "qwe"),
    @("5", "FOOD",      "44", "43"),
    @("6", "1ST FLOOR", "nad","nad")
)

$r = 3
foreach ($row in $data) {
    $rng = $ws.Range("A$r" + ":D$r")
    # Force text storage (values like "2", "44" would otherwise be
    # auto-coerced to numbers), then drop the formatting again so the
    # cells keep the workbook's default style.
    $rng.NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $rng.ClearFormats()
    $r++
}
